$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '34.457.24'
$ws.Cells.Item(2, 5).Value = '  +0.96%  '
$ws.Cells.Item(3, 4).Value = '1.796.41'
$ws.Cells.Item(3, 5).Value = '  +0.43%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '226.93'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.03%  '
$ws.Cells.Item(6, 5).Value = '  +1.56%  '
$ws.Cells.Item(7, 5).Value = '  +0.01%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '32.45'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +1.65%  '
$ws.Cells.Item(9, 5).Value = '  +1.34%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0694'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +0.73%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0949'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.38%  '
$ws.Cells.Item(12, 4).Value = '2.054.56'
$ws.Cells.Item(12, 5).Value = '  +0.38%  '
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = '1.826.02'
$ws.Cells.Item(13, 5).Value = '  +2.08%  '
$ws.Cells.Item(14, 2).Value = 'Chainlink'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '11.10'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -0.96%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.634'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +2.10%  '
$ws.Cells.Item(16, 4).Value = '34.414.07'
$ws.Cells.Item(16, 5).Value = '  +1.04%  '
$ws.Cells.Item(17, 5).Value = '  +1.58%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '68.41'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +0.41%  '
$ws.Cells.Item(19, 2).Value = 'BitcoinCash'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '246.66'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.39%  '
$ws.Cells.Item(20, 2).Value = 'ShibaInu'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(20, 4).Value = '0.0₃0800'
$ws.Cells.Item(20, 5).Value = '  +3.01%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '11.13'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +2.26%  '
$ws.Cells.Item(22, 5).Value = '  -0.02%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.17'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +1.77%  '
$ws.Cells.Item(24, 5).Value = '  +1.19%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '162.81'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.79%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.26'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +1.35%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '16.45'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +0.75%  '
$ws.Cells.Item(28, 5).Value = '  +2.28%  '
$ws.Cells.Item(29, 5).Value = '  +0.05%  '
$ws.Cells.Item(30, 5).Value = '  +0.21%  '
$ws.Cells.Item(32, 5).Value = '  +8.17%  '
$ws.Cells.Item(33, 5).Value = '  +3.02%  '
$ws.Cells.Item(34, 5).Value = '  +1.17%  '
$ws.Cells.Item(35, 4).Value = '1.442.95'
$ws.Cells.Item(35, 5).Value = '  -1.30%  '
$ws.Cells.Item(36, 5).Value = '  +6.70%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.670'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +3.74%  '
$ws.Cells.Item(38, 5).Value = '  -0.38%  '
$ws.Cells.Item(39, 5).Value = '  +1.76%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '84.04'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +4.57%  '
$ws.Cells.Item(41, 5).Value = '  +1.40%  '
$ws.Cells.Item(42, 5).Value = '  +1.54%  '
$ws.Cells.Item(43, 5).Value = '  +2.51%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '13.77'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +2.11%  '
$ws.Cells.Item(45, 5).Value = '  +3.62%  '
$ws.Cells.Item(46, 5).Value = '  +0.69%  '
$ws.Cells.Item(47, 5).Value = '  +0.14%  '
$ws.Cells.Item(48, 4).Value = '1.950.32'
$ws.Cells.Item(48, 5).Value = '  +0.11%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '105.73'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -0.57%  '
$ws.Cells.Item(50, 5).Value = '  -0.01%  '
$ws.Cells.Item(51, 4).Value = '0.0₆0130'
$ws.Cells.Item(51, 5).Value = '  -4.68%  '

Write-Host "Applied changes"
